$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D6").Value = "U"
$ws.Range("E2:G6").ClearContents()
$ws.Range("E2:E6").Value = 1

$ws.Range("L7").Select()
